$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2315
$ws.Range("E2").Value = 141
$ws.Range("F2").Value = 141
$ws.Range("G2").Value = 139
$ws.Range("H2").Value = 113
$ws.Range("I2").Value = 113
$ws.Range("K2").Value = 2925
$ws.Range("L2").Value = 642
$ws.Range("M2").Value = 2284
$ws.Range("N2").Value = 2284
$ws.Range("P2").Value = 200
$ws.Range("Q2").Value = 198
$ws.Range("R2").Value = -260
$ws.Range("S2").Value = 23
$ws.Range("T2").Value = 260
$ws.Range("U2").Value = -62
$ws.Range("V2").Value = 272
$ws.Range("W2").Value = 6.11
$ws.Range("X2").Value = 4.86
$ws.Range("Y2").Value = 5.03
$ws.Range("Z2").Value = 3.93
$ws.Range("AA2").Value = 28.09
$ws.Range("AB2").Value = 1041.77
$ws.Range("AC2").Value = 2814
$ws.Range("AD2").Value = 8.35
$ws.Range("AE2").Value = 57088
$ws.Range("AF2").Value = 0.41
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 2.13
$ws.Range("AI2").Value = 17.77
$ws.Range("AJ2").Value = 4000000
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 2307
$ws.Range("E3").Value = 27
$ws.Range("F3").Value = 27
$ws.Range("G3").Value = -4
$ws.Range("H3").Value = -2
$ws.Range("I3").Value = -2
$ws.Range("K3").Value = 3010
$ws.Range("L3").Value = 751
$ws.Range("M3").Value = 2259
$ws.Range("N3").Value = 2259
$ws.Range("P3").Value = 200
$ws.Range("Q3").Value = 139
$ws.Range("R3").Value = -84
$ws.Range("S3").Value = -35
$ws.Range("T3").Value = 94
$ws.Range("U3").Value = 44
$ws.Range("V3").Value = 258
$ws.Range("W3").Value = 1.17
$ws.Range("X3").Value = -0.1
$ws.Range("Y3").Value = -0.11
$ws.Range("Z3").Value = -0.08
$ws.Range("AA3").Value = 33.27
$ws.Range("AB3").Value = 1029.49
$ws.Range("AC3").Value = -60
$ws.Range("AD3").Value = -332.44
$ws.Range("AE3").Value = 56475
$ws.Range("AF3").Value = 0.35
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 2.51
$ws.Range("AI3").Value = -833.19
$ws.Range("AJ3").Value = 4000000
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 2283
$ws.Range("E4").Value = 31
$ws.Range("F4").Value = 31
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = -2
$ws.Range("I4").Value = -2
$ws.Range("K4").Value = 2997
$ws.Range("L4").Value = 756
$ws.Range("M4").Value = 2241
$ws.Range("N4").Value = 2241
$ws.Range("P4").Value = 200
$ws.Range("Q4").Value = 136
$ws.Range("R4").Value = -107
$ws.Range("S4").Value = -25
$ws.Range("T4").Value = 33
$ws.Range("U4").Value = 103
$ws.Range("V4").Value = 253
$ws.Range("W4").Value = 1.34
$ws.Range("X4").Value = -0.09
$ws.Range("Y4").Value = -0.09
$ws.Range("Z4").Value = -0.07
$ws.Range("AA4").Value = 33.72
$ws.Range("AB4").Value = 1018.75
$ws.Range("AC4").Value = -52
$ws.Range("AD4").Value = -345.03
$ws.Range("AE4").Value = 56026
$ws.Range("AF4").Value = 0.32
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 2.79
$ws.Range("AI4").Value = -963.77
$ws.Range("AJ4").Value = 4000000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 2444
$ws.Range("E5").Value = -69
$ws.Range("F5").Value = -69
$ws.Range("G5").Value = -67
$ws.Range("H5").Value = -53
$ws.Range("I5").Value = -53
$ws.Range("K5").Value = 2991
$ws.Range("L5").Value = 822
$ws.Range("M5").Value = 2169
$ws.Range("N5").Value = 2169
$ws.Range("P5").Value = 200
$ws.Range("Q5").Value = -78
$ws.Range("R5").Value = -71
$ws.Range("S5").Value = 119
$ws.Range("T5").Value = 55
$ws.Range("U5").Value = -132
$ws.Range("V5").Value = 391
$ws.Range("W5").Value = -2.82
$ws.Range("X5").Value = -2.16
$ws.Range("Y5").Value = -2.39
$ws.Range("Z5").Value = -1.76
$ws.Range("AA5").Value = 37.89
$ws.Range("AB5").Value = 983.62
$ws.Range("AC5").Value = -1318
$ws.Range("AD5").Value = -11.87
$ws.Range("AE5").Value = 54230
$ws.Range("AF5").Value = 0.29
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 3.19
$ws.Range("AI5").Value = -37.93
$ws.Range("AJ5").Value = 4000000
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 2802
$ws.Range("E6").Value = 231
$ws.Range("F6").Value = 231
$ws.Range("G6").Value = 231
$ws.Range("H6").Value = 173
$ws.Range("I6").Value = 173
$ws.Range("K6").Value = 2988
$ws.Range("L6").Value = 669
$ws.Range("M6").Value = 2320
$ws.Range("N6").Value = 2320
$ws.Range("P6").Value = 200
$ws.Range("Q6").Value = 420
$ws.Range("R6").Value = -5
$ws.Range("S6").Value = -292
$ws.Range("T6").Value = 20
$ws.Range("U6").Value = 400
$ws.Range("V6").Value = 119
$ws.Range("W6").Value = 8.23
$ws.Range("X6").Value = 6.17
$ws.Range("Y6").Value = 7.7
$ws.Range("Z6").Value = 5.78
$ws.Range("AA6").Value = 28.82
$ws.Range("AB6").Value = 1059.77
$ws.Range("AC6").Value = 4323
$ws.Range("AD6").Value = 4.5
$ws.Range("AE6").Value = 57993
$ws.Range("AF6").Value = 0.34
$ws.Range("AG6").Value = 600
$ws.Range("AH6").Value = 3.08
$ws.Range("AI6").Value = 13.88
$ws.Range("AJ6").Value = 4000000
$ws.Range("J6").ClearContents()
$ws.Range("O6").ClearContents()

# Rows 7-9: remove all data beyond column C (financial figures no longer reported)
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
